$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6148512105072825
$ws.Range("C2").Value = 0.1507818367487914
$ws.Range("D2").Value = 0.007246546123949926
$ws.Range("F2").Value = 0.7263384973622777
$ws.Range("G2").Value = 0.5877646244443184
$ws.Range("H2").Value = 0.6321235378771206
$ws.Range("M2").Value = 0.8925796379877795
$ws.Range("B3").Value = 0.5394337338853461
$ws.Range("C3").Value = 0.1351597823191923
$ws.Range("D3").Value = 0.007173905759168875
$ws.Range("F3").Value = 0.7021404881570845
$ws.Range("G3").Value = 0.5622255763363029
$ws.Range("H3").Value = 0.6254983670851715
$ws.Range("M3").Value = 0.7942143616790389
$ws.Range("B4").Value = 0.4931110360349464
$ws.Range("C4").Value = 0.1254945084185124
$ws.Range("D4").Value = 0.007129622630934307
$ws.Range("F4").Value = 0.6878717134335375
$ws.Range("G4").Value = 0.5470701036518761
$ws.Range("H4").Value = 0.6218670413900043
$ws.Range("M4").Value = 0.7343662059867597
$ws.Range("B5").Value = 0.474230932035141
$ws.Range("C5").Value = 0.1215376569476803
$ws.Range("D5").Value = 0.007111667324453563
$ws.Range("F5").Value = 0.6822039177746149
$ws.Range("G5").Value = 0.5410247882208949
$ws.Range("H5").Value = 0.6204964398665709
$ws.Range("M5").Value = 0.7101088043355759
$ws.Range("B6").Value = 0.4710957343623932
$ws.Range("C6").Value = 0.1208795345456224
$ws.Range("D6").Value = 0.007108691633787245
$ws.Range("F6").Value = 0.6812716176143141
$ws.Range("G6").Value = 0.5400288159435291
$ws.Range("H6").Value = 0.6202754305314926
$ws.Range("M6").Value = 0.7060886119116958
$ws.Range("B7").Value = 0.4928564243364235
$ws.Range("C7").Value = 0.1254412181885982
$ws.Range("D7").Value = 0.007129380098508165
$ws.Range("F7").Value = 0.687794682616655
$ws.Range("G7").Value = 0.546988047350581
$ws.Range("H7").Value = 0.6218481156248856
$ws.Range("M7").Value = 0.7340385397209275
$ws.Range("B8").Value = 0.5888510565893625
$ws.Range("C8").Value = 0.1454107041364239
$ws.Range("D8").Value = 0.007221439714875544
$ws.Range("F8").Value = 0.7178720427961593
$ws.Range("G8").Value = 0.578848755981042
$ws.Range("H8").Value = 0.629748217223451
$ws.Range("M8").Value = 0.8585453317533336
$ws.Range("B9").Value = 0.7769431924256196
$ws.Range("C9").Value = 0.1839803221987779
$ws.Range("D9").Value = 0.007404068904229177
$ws.Range("F9").Value = 0.7815859350519929
$ws.Range("G9").Value = 0.6455694448135603
$ws.Range("H9").Value = 0.6487320463130573
$ws.Range("M9").Value = 1.10737509715689
$ws.Range("B10").Value = 0.9150215292424377
$ws.Range("C10").Value = 0.2119479289073922
$ws.Range("D10").Value = 0.007538999927181322
$ws.Range("F10").Value = 0.8313694635904483
$ws.Range("G10").Value = 0.6972793110291207
$ws.Range("H10").Value = 0.6648490679119448
$ws.Range("M10").Value = 1.293506760867913
$ws.Range("B11").Value = 0.9778094000730562
$ws.Range("C11").Value = 0.2245890779925048
$ws.Range("D11").Value = 0.007600445563582525
$ws.Range("F11").Value = 0.8546821875950599
$ws.Range("G11").Value = 0.7214104444005329
$ws.Range("H11").Value = 0.6726612419012383
$ws.Range("M11").Value = 1.37900828285386
$ws.Range("B12").Value = 1.001581438519111
$ws.Range("C12").Value = 0.2293640320962993
$ws.Range("D12").Value = 0.007623714702337026
$ws.Range("F12").Value = 0.8636073365930628
$ws.Range("G12").Value = 0.7306374645151266
$ws.Range("H12").Value = 0.6756892905970346
$ws.Range("M12").Value = 1.411513248162692
$ws.Range("B13").Value = 0.9964619117529878
$ws.Range("C13").Value = 0.228336197006513
$ws.Range("D13").Value = 0.007618703359412393
$ws.Range("F13").Value = 0.8616808045237008
$ws.Range("G13").Value = 0.7286462723030809
$ws.Range("H13").Value = 0.6750340335305793
$ws.Range("M13").Value = 1.404506924319733
$ws.Range("B14").Value = 0.9797652321169608
$ws.Range("C14").Value = 0.2249821578192268
$ws.Range("D14").Value = 0.007602359955203752
$ws.Range("F14").Value = 0.8554145105720181
$ws.Range("G14").Value = 0.7221677620656237
$ws.Range("H14").Value = 0.6729089596481117
$ws.Range("M14").Value = 1.381679885278785
$ws.Range("B15").Value = 0.9695374411380158
$ws.Range("C15").Value = 0.2229261447980662
$ws.Range("D15").Value = 0.007592349041914304
$ws.Range("F15").Value = 0.8515889139495414
$ws.Range("G15").Value = 0.7182111378235732
$ws.Range("H15").Value = 0.6716163934049746
$ws.Range("M15").Value = 1.367714500842069
$ws.Range("B16").Value = 0.9109176276140829
$ws.Range("C16").Value = 0.211120138571772
$ws.Range("D16").Value = 0.007534984893631247
$ws.Range("F16").Value = 0.8298594391993248
$ws.Range("G16").Value = 0.6957146608185667
$ws.Range("H16").Value = 0.664348248215731
$ws.Range("M16").Value = 1.287936413814322
$ws.Range("B17").Value = 0.8749493671323876
$ws.Range("C17").Value = 0.2038564864620582
$ws.Range("D17").Value = 0.007499804385862063
$ws.Range("F17").Value = 0.8167006429484616
$ws.Range("G17").Value = 0.6820707064721887
$ws.Range("H17").Value = 0.6600130007196583
$ws.Range("M17").Value = 1.239213266832053
$ws.Range("B18").Value = 0.8542591021897579
$ws.Range("C18").Value = 0.1996709743495728
$ws.Range("D18").Value = 0.007479576057416892
$ws.Range("F18").Value = 0.8091946514182808
$ws.Range("G18").Value = 0.6742801989161649
$ws.Range("H18").Value = 0.657564664540871
$ws.Range("M18").Value = 1.211266557967676
$ws.Range("B19").Value = 0.8472533688648696
$ws.Range("C19").Value = 0.198252525722296
$ws.Range("D19").Value = 0.007472728482948554
$ws.Range("F19").Value = 0.8066639618547384
$ws.Range("G19").Value = 0.6716522297653853
$ws.Range("H19").Value = 0.6567434403077357
$ws.Range("M19").Value = 1.201817358531656
$ws.Range("B20").Value = 0.8787784910808227
$ws.Range("C20").Value = 0.2046305084363382
$ws.Range("D20").Value = 0.007503548792641368
$ws.Range("F20").Value = 0.8180949305574359
$ws.Range("G20").Value = 0.6835172036511494
$ws.Range("H20").Value = 0.6604698141427434
$ws.Range("M20").Value = 1.244391843919942
$ws.Range("B21").Value = 0.9846695765971276
$ws.Range("C21").Value = 0.2259676474920411
$ws.Range("D21").Value = 0.007607160445889605
$ws.Range("F21").Value = 0.8572524261045231
$ws.Range("G21").Value = 0.7240682266913154
$ws.Range("H21").Value = 0.6735312471466557
$ws.Range("M21").Value = 1.388381216509941
$ws.Range("B22").Value = 1.053849706881579
$ws.Range("C22").Value = 0.2398427841956163
$ws.Range("D22").Value = 0.007674879014828662
$ws.Range("F22").Value = 0.8834108796017404
$ws.Range("G22").Value = 0.7510906832899877
$ws.Range("H22").Value = 0.6824745024631511
$ws.Range("M22").Value = 1.483233090979823
$ws.Range("B23").Value = 1.016929598416482
$ws.Range("C23").Value = 0.2324438457050064
$ws.Range("D23").Value = 0.007638738746685192
$ws.Range("F23").Value = 0.8693973083705941
$ws.Range("G23").Value = 0.7366201613218095
$ws.Range("H23").Value = 0.6776638727093598
$ws.Range("M23").Value = 1.432537806068098
$ws.Range("B24").Value = 0.8770473803151049
$ws.Range("C24").Value = 0.2042806028024984
$ws.Range("D24").Value = 0.007501855954256342
$ws.Range("F24").Value = 0.8174643889132938
$ws.Range("G24").Value = 0.6828630755001086
$ws.Range("H24").Value = 0.6602631516211943
$ws.Range("M24").Value = 1.242050407753183
$ws.Range("B25").Value = 0.7260777650667478
$ws.Range("C25").Value = 0.1736104302616468
$ws.Range("D25").Value = 0.007354494293874581
$ws.Range("F25").Value = 0.7638339452725944
$ws.Range("G25").Value = 0.6270548626269772
$ws.Range("H25").Value = 0.6432183852403455
$ws.Range("M25").Value = 1.039512179900314
